# Auto commit at 2025-09-12  8:15:45.50
# Refresh the daily "Metrics" figures (Sheet: Metrics) with the latest
# cumulative totals, then leave the selection where the operator's cursor
# landed on each sheet. Everything downstream (the "today" sheet's
# Metrics! lookups and running totals, and TODAY()-1) recalculates on its
# own once the source cells change.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Activate()

$metrics.Range("B2").Value2  = 171717.73
$metrics.Range("B3").Value2  = 139048.24000000002
$metrics.Range("B4").Value2  = 54491.520000000004
$metrics.Range("B5").Value2  = 6742
$metrics.Range("B6").Value2  = 4090968.61
$metrics.Range("B7").Value2  = 3466575.7199999993
$metrics.Range("B8").Value2  = 1183857.2
$metrics.Range("B9").Value2  = 157902
$metrics.Range("B10").Value2 = 32556292.410999827
$metrics.Range("B11").Value2 = 19496445.790000003
$metrics.Range("B12").Value2 = 11465566.090000002
$metrics.Range("B13").Value2 = 1255529

$metrics.Range("F17").Select() | Out-Null

$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("E6").Select() | Out-Null
